# "New forms for registration"
# Updates the "choices" sheet with new delivery-site choices (Ritsona, Skaramagas,
# Nea Kavala, Softex, Athens replace the old a0..a3 / Site A..D codes), marks the
# "delivery_site" and a new question as required on the survey sheet, and widens a
# couple of columns / updates selections accordingly.

$wb = $excel.ActiveWorkbook

# --- survey sheet -----------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Column K ("required") gets a dedicated width and the header takes the same
# style as the other header cells (H1/I1/J1 use style index 10).
$survey.Columns.Item(11).ColumnWidth = 21.33203125
$survey.Range("K1").Style = $survey.Range("J1").Style

# Mark "telephone" (row 13) and the new "tent_caravan" (row 22) rows as required.
$survey.Range("K13").Value = $true
$survey.Range("K22").Value = $true

# Move the current selection / top row, as left by the author.
$survey.Application.ActiveWindow.ScrollRow = 7
$survey.Range("E13").Select()

# --- choices sheet ------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("B5").Value = "Ritsona"
$choices.Range("C5").Value = "Ritsona"

$choices.Range("B6").Value = "Skaramagas"
$choices.Range("C6").Value = "Skaramagas"

$choices.Range("B7").Value = "Nea Kavala"
$choices.Range("C7").Value = "Nea Kavala"

$choices.Range("B8").Value = "Softex"
$choices.Range("C8").Value = "Softex"

$choices.Range("A9").Value = "sites"
$choices.Range("B9").Value = "Athens"
$choices.Range("C9").Value = "Athens"

$choices.Range("C7").Select()
